$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 51226.45
$ws.Range("I6").Value = 157
$ws.Range("K6").Value = 471
$ws.Range("M6").Value = -359

# Row 40
$ws.Range("H40").Value = 5000
$ws.Range("I40").Value = 4000
$ws.Range("K40").Value = 4000
$ws.Range("M40").Value = -3825

# Row 64
$ws.Range("H64").Value = 40280.438
$ws.Range("I64").Value = 48768.23
$ws.Range("J64").Value = 3500
$ws.Range("K64").Value = 48768.23
$ws.Range("L64").Value = 3500
$ws.Range("M64").Value = -48520.23
$ws.Range("N64").Value = -3996

# Row 67
$ws.Range("H67").Value = 40280.438
$ws.Range("I67").Value = 48768.23
$ws.Range("J67").Value = 3500
$ws.Range("K67").Value = 48768.23
$ws.Range("L67").Value = 3500
$ws.Range("M67").Value = -47910.23
$ws.Range("N67").Value = -5216

# Row 101
$ws.Range("H101").Value = 20409944
$ws.Range("I101").Value = 23811436
$ws.Range("J101").Value = 999
$ws.Range("K101").Value = 71434308
$ws.Range("L101").Value = 2997
$ws.Range("M101").Value = -71432686
$ws.Range("N101").Value = -6241

# Row 121
$ws.Range("H121").Value = 1618.5454
$ws.Range("J121").Value = 1618.5454
$ws.Range("L121").Value = 4855.6362
$ws.Range("N121").Value = -8349.636200000001

# Row 132
$ws.Range("H132").Value = 1727223.6
$ws.Range("I132").Value = 2758.8333
$ws.Range("K132").Value = 8276.499899999999
$ws.Range("M132").Value = -5746.499899999999

# Row 137
$ws.Range("I137").Value = 1022215.56
$ws.Range("J137").Value = 3811.2222
$ws.Range("K137").Value = 3066646.68
$ws.Range("L137").Value = 11433.6666
$ws.Range("M137").Value = -3064096.68
$ws.Range("N137").Value = -16533.6666

# Row 138
$ws.Range("H138").Value = 195141.25
$ws.Range("I138").Value = 574337.25
$ws.Range("J138").Value = 5543.2393
$ws.Range("K138").Value = 1723011.75
$ws.Range("L138").Value = 16629.7179
$ws.Range("M138").Value = -1717871.75
$ws.Range("N138").Value = -26909.7179

# Row 141
$ws.Range("H141").Value = 7003.8
$ws.Range("I141").Value = 5942.25
$ws.Range("K141").Value = 17826.75
$ws.Range("M141").Value = -12646.75


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9496.606
$ws.Range("I32").Value = 8699.636
$ws.Range("K32").Value = 8699.636
$ws.Range("M32").Value = -8412.636

# Row 96
$ws.Range("H96").Value = 62829.668
$ws.Range("J96").Value = 62829.668
$ws.Range("L96").Value = 62829.668
$ws.Range("N96").Value = -68321.66800000001

# Row 110
$ws.Range("H110").Value = 2396.3845
$ws.Range("I110").Value = 1566.25
$ws.Range("J110").Value = 3724.6
$ws.Range("K110").Value = 1566.25
$ws.Range("L110").Value = 3724.6
$ws.Range("M110").Value = 478.75
$ws.Range("N110").Value = -7814.6

# Row 124
$ws.Range("H124").Value = 50000
$ws.Range("J124").Value = 50000
$ws.Range("L124").Value = 50000
$ws.Range("N124").Value = -59820

# Row 132
$ws.Range("H132").Value = 2183.7036
$ws.Range("I132").Value = 1873.15
$ws.Range("J132").Value = 3071
$ws.Range("K132").Value = 5619.450000000001
$ws.Range("L132").Value = 9213
$ws.Range("M132").Value = -3089.450000000001
$ws.Range("N132").Value = -14273


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 42
$ws.Range("H42").Value = 194760.8
$ws.Range("J42").Value = 194760.8
$ws.Range("L42").Value = 194760.8
$ws.Range("N42").Value = -195416.8

# Row 134
$ws.Range("H134").Value = 3207.6897
$ws.Range("I134").Value = 2733.0908
$ws.Range("K134").Value = 8199.2724
$ws.Range("M134").Value = -5664.2724


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2687.125
$ws.Range("I31").Value = 1455
$ws.Range("J31").Value = 5397.8
$ws.Range("K31").Value = 1455
$ws.Range("L31").Value = 5397.8
$ws.Range("M31").Value = -1160
$ws.Range("N31").Value = -5987.8

# Row 34
$ws.Range("H34").Value = 2687.125
$ws.Range("I34").Value = 1455
$ws.Range("J34").Value = 5397.8
$ws.Range("K34").Value = 1455
$ws.Range("L34").Value = 5397.8
$ws.Range("M34").Value = -1253
$ws.Range("N34").Value = -5801.8

# Row 41
$ws.Range("H41").Value = 18666.666
$ws.Range("J41").Value = 23000
$ws.Range("L41").Value = 23000
$ws.Range("N41").Value = -23856

# Row 122
$ws.Range("H122").Value = 6936.613
$ws.Range("I122").Value = 8542.299999999999
$ws.Range("J122").Value = 4017.182
$ws.Range("K122").Value = 25626.9
$ws.Range("L122").Value = 12051.546
$ws.Range("M122").Value = -23176.9
$ws.Range("N122").Value = -16951.546

# Row 134
$ws.Range("H134").Value = 3314.4614
$ws.Range("I134").Value = 3382.3333
$ws.Range("K134").Value = 10146.9999
$ws.Range("M134").Value = -7611.999899999999

# Row 141
$ws.Range("H141").Value = 447944.44
$ws.Range("J141").Value = 592087.8
$ws.Range("L141").Value = 592087.8
$ws.Range("N141").Value = -602447.8


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 716558.4399999999
$ws.Range("I5").Value = 1256.5
$ws.Range("J5").Value = 1253034.9
$ws.Range("K5").Value = 3769.5
$ws.Range("L5").Value = 3759104.7
$ws.Range("M5").Value = -3657.5
$ws.Range("N5").Value = -3759328.7

# Row 33
$ws.Range("H33").Value = 60.666668
$ws.Range("I33").Value = 33.5
$ws.Range("K33").Value = 201
$ws.Range("M33").Value = 82

# Row 49
$ws.Range("H49").Value = 200
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()

# Row 122
$ws.Range("H122").Value = 2600.0286
$ws.Range("J122").Value = 3127.5
$ws.Range("L122").Value = 28147.5
$ws.Range("N122").Value = -33047.5

# Row 129
$ws.Range("H129").Value = 1510.6154
$ws.Range("I129").Value = 1239.909
$ws.Range("J129").Value = 2999.5
$ws.Range("K129").Value = 3719.727
$ws.Range("L129").Value = 8998.5
$ws.Range("M129").Value = 1280.273
$ws.Range("N129").Value = -18998.5

# Row 132
$ws.Range("H132").Value = 70355.664
$ws.Range("I132").Value = 2000
$ws.Range("J132").Value = 78900.125
$ws.Range("K132").Value = 18000
$ws.Range("L132").Value = 710101.125
$ws.Range("M132").Value = -15470
$ws.Range("N132").Value = -715161.125

# Row 135
$ws.Range("H135").Value = 716558.4399999999
$ws.Range("I135").Value = 1256.5
$ws.Range("J135").Value = 1253034.9
$ws.Range("K135").Value = 11308.5
$ws.Range("L135").Value = 11277314.1
$ws.Range("M135").Value = -8773.5
$ws.Range("N135").Value = -11282384.1


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 10000
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

# Row 97
$ws.Range("H97").Value = 9399.462
$ws.Range("I97").Value = 10841.091
$ws.Range("K97").Value = 10841.091
$ws.Range("M97").Value = -10345.091

# Row 132
$ws.Range("H132").Value = 4097.5
$ws.Range("I132").Value = 4238.3706
$ws.Range("J132").Value = 3674.889
$ws.Range("K132").Value = 12715.1118
$ws.Range("L132").Value = 11024.667
$ws.Range("M132").Value = -10185.1118
$ws.Range("N132").Value = -16084.667


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 12596.538
$ws.Range("I16").Value = 13886.909
$ws.Range("J16").Value = 5499.5
$ws.Range("K16").Value = 13886.909
$ws.Range("L16").Value = 5499.5
$ws.Range("M16").Value = -13716.909
$ws.Range("N16").Value = -5839.5

# Row 55
$ws.Range("H55").Value = 4931.125
$ws.Range("I55").Value = 949.25
$ws.Range("K55").Value = 949.25
$ws.Range("M55").Value = -776.25

# Row 61
$ws.Range("H61").Value = 17803.588
$ws.Range("I61").Value = 2276.25
$ws.Range("J61").Value = 31605.666
$ws.Range("K61").Value = 2276.25
$ws.Range("L61").Value = 31605.666
$ws.Range("M61").Value = -2074.25
$ws.Range("N61").Value = -32009.666

# Row 100
$ws.Range("H100").Value = 4034.6365
$ws.Range("I100").Value = 1422.625
$ws.Range("K100").Value = 1422.625
$ws.Range("M100").Value = -881.625

# Row 113
$ws.Range("H113").Value = 17803.588
$ws.Range("I113").Value = 2276.25
$ws.Range("J113").Value = 31605.666
$ws.Range("K113").Value = 2276.25
$ws.Range("L113").Value = 31605.666
$ws.Range("M113").Value = -106.25
$ws.Range("N113").Value = -35945.666

# Row 122
$ws.Range("H122").Value = 5998.5557
$ws.Range("I122").Value = 5711.8184
$ws.Range("J122").Value = 6449.143
$ws.Range("K122").Value = 17135.4552
$ws.Range("L122").Value = 19347.429
$ws.Range("M122").Value = -14685.4552
$ws.Range("N122").Value = -24247.429

# Row 132
$ws.Range("H132").Value = 537559.1
$ws.Range("I132").Value = 1248256.9
$ws.Range("K132").Value = 3744770.7
$ws.Range("M132").Value = -3742240.7

# Row 136
$ws.Range("H136").Value = 7366.515
$ws.Range("I136").Value = 5294.5713
$ws.Range("J136").Value = 10992.417
$ws.Range("K136").Value = 15883.7139
$ws.Range("L136").Value = 32977.251
$ws.Range("M136").Value = -13333.7139
$ws.Range("N136").Value = -38077.251


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 52
$ws.Range("H52").Value = 25500
$ws.Range("I52").Value = 25500
$ws.Range("K52").Value = 25500
$ws.Range("M52").Value = -25274

# Row 96
$ws.Range("H96").Value = 2278.5715
$ws.Range("I96").Value = 2354.5454
$ws.Range("K96").Value = 2354.5454
$ws.Range("M96").Value = -981.5454

# Row 136
$ws.Range("H136").Value = 735691.9399999999
$ws.Range("I136").Value = 857307.25
$ws.Range("K136").Value = 2571921.75
$ws.Range("M136").Value = -2569371.75

